$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for the team-record columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/bordered/centered header style used by the rest of row 1
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row shares the same team record for this season
for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 30).Value = 80
    $ws.Cells.Item($row, 31).Value = 82
    $ws.Cells.Item($row, 32).Value = 0
}
